$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.967.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "'2.356.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'545.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").Value = "'132.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "'2.352.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "'0.334"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "'23.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "'2.768.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "'59.952.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "'2.351.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").Value = "'10.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "'4.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").Value = "'6.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.32%  "
$ws.Range("D22").Value = "'313.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'63.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "'0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'7.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("E28").Value = "  +5.05%  "
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("D30").Value = "'171.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.40%  "
$ws.Range("D32").Value = "'0.0₃0728"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'5.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.22%  "
$ws.Range("D34").Value = "'1.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.09%  "
$ws.Range("D35").Value = "'0.382"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "'18.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +6.74%  "
$ws.Range("D40").Value = "'314.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.84%  "
$ws.Range("D41").Value = "'38.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").Value = "'142.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'3.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "'0.0951"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "'19.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.03%  "
$ws.Range("D47").Value = "'0.0498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'0.563"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "'11.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "'0.0₆0209"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.97%  "
